# Append a new data row (row 4) to the Site Total Alarms sheet, matching
# the existing rows' layout: col A is a text-formatted date, columns
# B:G are plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: text date, styled exactly like A2/A3 -----------------
# Clone A2's cell style (border/font/alignment) onto A4 first.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats

# Force the incoming value to be stored as text (not auto-converted to
# a number) the same way the source rows are stored.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "20240223"

# Re-apply the original (General) cell format on top so the cell ends
# up sharing the very same style as A2/A3 rather than a bespoke "Text"
# number format.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats

# --- Columns B:G: plain numbers --------------------------------------
$ws.Range("B4").Value = 174
$ws.Range("C4").Value = 674
$ws.Range("D4").Value = 345
$ws.Range("E4").Value = 86
$ws.Range("F4").Value = 45
$ws.Range("G4").Value = 152
